# edit.ps1 - apply "Work tracker" worklog update
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace the 8 "Work tracker" bullet paragraphs (list item 1..8)
#    with the new 10-paragraph block described by the commit:
#      - para 1 becomes a plain (non-bulleted) paragraph holding the
#        original "Work tracker ... INDEX.HTML" text
#      - the remaining bullets are re-worded / re-ordered / extended
#      - two brand new bullets are appended at the end
# ------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$lastPara  = $d.Paragraphs.Item(8)
$blockRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)

$newBlockXml = '<w:p><w:r><w:t xml:space="preserve">Work tracker </w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>Ash</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>Declan</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>Miller</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Speedie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t xml:space="preserve">Brad </w:t></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>Facilitated creation of INDEX.HTML</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Extracted Bio/About me pages from team </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GitHubs</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Processed / formatted </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GitHubs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (BIO)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t>Truncated GitHub profiles (BIO)</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Utilised personal Grammarly account to clean up grammar / language </w:t></w:r><w:r><w:br/><w:t>Uploaded source data to web site</w:t></w:r><w:r><w:br/><w:t>Posted feedback on project feedback and voted</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">BIO information now </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>live</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> (Crude – needs work)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Processing Ash’s IT Security write up in Grammarly (approx. 150 positive aspects to update) </w:t></w:r><w:r w:rsidR="000327CF"><w:t>– this will take time to process</w:t></w:r><w:r><w:t xml:space="preserve"> – completed this</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Performed test for team profile to make </w:t></w:r><w:r><w:t xml:space="preserve">bio </w:t></w:r><w:r><w:t xml:space="preserve">available </w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Created readme.MD for GitHub</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Very actively </w:t></w:r><w:r><w:t xml:space="preserve">participated in </w:t></w:r><w:r><w:t>Trello</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Agreed to take on IT Tech – Cloud services, computing. Wrote 600 words to upload to web page – performed in </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>2 hour</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> study period and added citations. Requires further revision. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Uploaded further files to site. </w:t></w:r></w:p>'

[void]$blockRange.Delete()
$insertionPoint = $d.Range(0, 0)
[void]$insertionPoint.InsertXML($newBlockXml)

# ------------------------------------------------------------------
# 2) Add a <w:lastRenderedPageBreak/> before the screenshot <w:drawing>
# ------------------------------------------------------------------
$shape = $d.InlineShapes.Item(1)
$picRange = $d.Range($shape.Range.Start, $shape.Range.Start)
[void]$picRange.InsertXML('<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:lastRenderedPageBreak/></w:r>')

Write-Output "done"
